$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.086.39'
$ws.Range("E2").Value = '  -1.42%  '

$ws.Range("D3").Value = '1.860.96'
$ws.Range("E3").Value = '  +0.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9965'
$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.25'
$ws.Range("E5").Value = '  +0.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9977'
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4645'
$ws.Range("E7").Value = '  -1.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2793'
$ws.Range("E8").Value = '  +1.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06428'
$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.51'
$ws.Range("E10").Value = '  +4.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '100.94'
$ws.Range("E11").Value = '  +19.11%  '

$ws.Range("D12").Value = '1.845.29'
$ws.Range("E12").Value = '  -0.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07534'
$ws.Range("E13").Value = '  +1.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.008'
$ws.Range("E14").Value = '  -0.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6309'
$ws.Range("E15").Value = '  +0.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '304.82'
$ws.Range("E16").Value = '  +25.61%  '

$ws.Range("D17").Value = '30.080.84'
$ws.Range("E17").Value = '  -1.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9983'
$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.66'
$ws.Range("E19").Value = '  -0.31%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007405'
$ws.Range("E20").Value = '  +0.66%  '

$ws.Range("D21").Value = '2.083.84'
$ws.Range("E21").Value = '  -0.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9969'
$ws.Range("E22").Value = '  -0.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.023'
$ws.Range("E23").Value = '  +1.44%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.146'
$ws.Range("E24").Value = '  +2.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '165.00'
$ws.Range("E25").Value = '  +1.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.115'
$ws.Range("E26").Value = '  -1.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.50'
$ws.Range("E27").Value = '  +8.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.943'
$ws.Range("E28").Value = '  +2.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1082'
$ws.Range("E29").Value = '  +6.08%  '

$ws.Range("E30").Value = '  -3.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.014'
$ws.Range("E31").Value = '  -0.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.864'
$ws.Range("E32").Value = '  +0.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04886'
$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7387'
$ws.Range("E34").Value = '  +4.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.121'
$ws.Range("E35").Value = '  -1.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.716'
$ws.Range("E36").Value = '  +0.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01905'
$ws.Range("E37").Value = '  +0.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.645'
$ws.Range("E38").Value = '  -1.54%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.973'
$ws.Range("E39").Value = '  -0.29%  '

$ws.Range("B40").Value = 'Quant'
$ws.Range("C40").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '107.04'
$ws.Range("E40").Value = '  +1.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8620'
$ws.Range("E41").Value = '  -1.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9976'
$ws.Range("E42").Value = '  -0.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.739'
$ws.Range("E43").Value = '  +4.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4064'
$ws.Range("E44").Value = '  -0.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.26'
$ws.Range("E45").Value = '  +5.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.103'
$ws.Range("E46").Value = '  -1.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.076'
$ws.Range("E47").Value = '  +6.70%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1195'
$ws.Range("E48").Value = '  -1.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.01'
$ws.Range("E49").Value = '  +2.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05554'
$ws.Range("E50").Value = '  +0.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3733'
$ws.Range("E51").Value = '  +1.27%  '
